$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")
$ws.Activate()

$ws.Range("C1").Value = "Must Read"
$ws.Range("C2").Value = "Yes"
$ws.Range("C3").Value = "No"

$ws.Range("C1").Font.Bold = $true

$ws.Range("C4").Select()
